# edit.ps1
# Applies the OOXML diff to the document:
#  1. Paragraph 1 ("To-do task 5"):
#       - removes the centered alignment (w:jc) from the paragraph
#       - splits the single run into 5 runs, adding leading/trailing
#         whitespace runs and a new "Lec: ArrayList 1,2" heading fragment
#  2. The "Setting adapter for arraylist and addfriend() function code."
#     bullet paragraph: splits the single run into 3 runs, wrapping
#     "addfriend(" with w:proofErr gramStart/gramEnd markers.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Locate paragraph 1 ("To-do task 5") and replace its run content.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1Range = $d.Range($p1.Range.Start, $p1.Range.End)

$frag1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">                                           </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>To-do task 5</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">                       </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Lec: </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>ArrayList 1,2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$p1Range.InsertXML($frag1)

# Remove the centered paragraph alignment (w:jc) -> default/left alignment.
$p1 = $d.Paragraphs(1)
$p1.Alignment = 0

# ---------------------------------------------------------------------
# 2) Locate the "addfriend()" bullet paragraph and replace its run
#    content, inserting proofErr gramStart/gramEnd markers around
#    "addfriend(".
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.TrimEnd() -eq "Setting adapter for arraylist and addfriend() function code.") {
        $target = $cand
        break
    }
}

$frag2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Setting adapter for arraylist and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>addfriend(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>) function code.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

if ($target -ne $null) {
    $targetRange = $d.Range($target.Range.Start, $target.Range.End)
    $targetRange.InsertXML($frag2)
}
